$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "67.398.90"
$ws.Range("E2").Value = "  -1.11%  "
# Row 3: Ethereum
$ws.Range("D3").Value = "3.231.68"
$ws.Range("E3").Value = "  -1.33%  "
# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.01%  "
# Row 5: BNB
$ws.Range("D5").Value = "'578.35"
$ws.Range("E5").Value = "  -1.63%  "
# Row 6: Solana
$ws.Range("D6").Value = "'183.99"
$ws.Range("E6").Value = "  -1.57%  "
# Row 7: XRP
$ws.Range("D7").Value = "'0.612"
$ws.Range("E7").Value = "  +1.70%  "
# Row 8: USDC
$ws.Range("E8").Value = "  +0.02%  "
# Row 9: LidoStakedEther
$ws.Range("D9").Value = "3.228.88"
$ws.Range("E9").Value = "  -1.34%  "
# Row 10: Dogecoin
$ws.Range("E10").Value = "  -3.62%  "
# Row 12: Cardano
$ws.Range("E12").Value = "  -1.64%  "
# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.784.88"
$ws.Range("E13").Value = "  -1.41%  "
# Row 14: TRON
$ws.Range("E14").Value = "  +0.03%  "
# Row 15: Avalanche
$ws.Range("D15").Value = "'27.56"
$ws.Range("E15").Value = "  -3.96%  "
# Row 16: WrappedBTC
$ws.Range("D16").Value = "67.460.96"
$ws.Range("E16").Value = "  -1.02%  "
# Row 17: ShibaInu
$ws.Range("E17").Value = "  -2.23%  "
# Row 18: WrappedEther
$ws.Range("D18").Value = "3.228.29"
$ws.Range("E18").Value = "  -1.25%  "
# Row 19: Polkadot
$ws.Range("E19").Value = "  -2.01%  "
# Row 20: Chainlink
$ws.Range("E20").Value = "  -1.41%  "
# Row 21: BitcoinCash
$ws.Range("D21").Value = "'394.51"
$ws.Range("E21").Value = "  +3.40%  "
# Row 22: Uniswap
$ws.Range("E22").Value = "  -2.58%  "
# Row 23: Dai
$ws.Range("E23").Value = "  -0.08%  "
# Row 24: Litecoin
$ws.Range("D24").Value = "'71.27"
$ws.Range("E24").Value = "  -0.42%  "
# Row 25: Polygon
$ws.Range("D25").Value = "'0.515"
$ws.Range("E25").Value = "  +0.05%  "
# Row 26: PEPE
$ws.Range("E26").Value = "  -3.01%  "
# Row 27: Kaspa
$ws.Range("D27").Value = "'0.188"
$ws.Range("E27").Value = "  -0.12%  "
# Row 28: InternetComputer(DFINITY)
$ws.Range("D28").Value = "'9.55"
$ws.Range("E28").Value = "  -2.66%  "
# Row 29: Binance-PegBSC-USD
$ws.Range("E29").Value = "  +0.19%  "
# Row 30: PancakeSwap
$ws.Range("E30").Value = "  -2.32%  "
# Row 31: NEARProtocol
$ws.Range("E31").Value = "  -4.51%  "
# Row 32: EthereumClassic
$ws.Range("D32").Value = "'22.54"
$ws.Range("E32").Value = "  -1.70%  "
# Row 33: Aptos
$ws.Range("D33").Value = "'6.97"
$ws.Range("E33").Value = "  -2.92%  "
# Row 34: USDe (was Fetch.AI)
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "  +0.00%  "
# Row 35: Fetch.AI (was USDe)
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").Value = "'1.25"
$ws.Range("E35").Value = "  -2.54%  "
# Row 36: Monero
$ws.Range("D36").Value = "'161.37"
$ws.Range("E36").Value = "  -1.02%  "
# Row 37: ImmutableX
$ws.Range("E37").Value = "  -4.81%  "
# Row 38: Stacks
$ws.Range("E38").Value = "  +0.59%  "
# Row 39: EnergySwap
$ws.Range("D39").Value = "'26.37"
$ws.Range("E39").Value = "  -1.14%  "
# Row 40: Mantle
$ws.Range("E40").Value = "  -4.36%  "
# Row 41: Filecoin
$ws.Range("E41").Value = "  -1.34%  "
# Row 42: RenderToken
$ws.Range("D42").Value = "'6.48"
$ws.Range("E42").Value = "  -4.89%  "
# Row 43: dogwifhat
$ws.Range("E43").Value = "  -6.29%  "
# Row 44: Hedera
$ws.Range("D44").Value = "'0.0684"
$ws.Range("E44").Value = "  -0.92%  "
# Row 45: OKB
$ws.Range("D45").Value = "'40.48"
$ws.Range("E45").Value = "  -1.89%  "
# Row 46: Maker
$ws.Range("D46").Value = "2.604.96"
$ws.Range("E46").Value = "  -1.67%  "
# Row 47: InjectiveProtocol
$ws.Range("D47").Value = "'24.63"
$ws.Range("E47").Value = "  -3.41%  "
# Row 48: Bittensor
$ws.Range("D48").Value = "'334.15"
$ws.Range("E48").Value = "  -3.91%  "
# Row 49: VeChain
$ws.Range("E49").Value = "  -2.41%  "
# Row 50: Stellar (was Cosmos)
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.102"
$ws.Range("E50").Value = "  -0.61%  "
# Row 51: Cosmos (was Stellar)
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "'6.26"
$ws.Range("E51").Value = "  -0.31%  "
